$d = $word.ActiveDocument

$subjectsText = "Subjects: (Subject, Occurrence : SubjectKind, Attribute : Resource P, Value : Resource O);"

# --- 1. Turn the "Subjects:" bullet into a new "Template Transforms:" header bullet,
#        followed (in the same paragraph) by a trailing empty run, and move the
#        original "Subjects:" text into a brand-new bullet paragraph right after it. ---

# 1a. Replace the text in-place (single run -> single run, clean replace).
$d.Content.Find.Execute($subjectsText, $true, $false, $false, $false, $false, $true, 1, $false, "Template Transforms:", 2)

# 1b. Insert two fresh empty paragraphs (same list formatting) right after it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Template Transforms:*") {
        $p.Range.InsertParagraphAfter()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Template Transforms:*") {
        $p2 = $p.Next()
        $p2.Range.InsertParagraphAfter()
        break
    }
}

# 1c. Merge the first new (empty) paragraph back into the "Template Transforms:"
#     paragraph by deleting the paragraph mark between them - this leaves the
#     "Template Transforms:" paragraph with its text run plus a second, empty
#     trailing run, and keeps the second new paragraph as a standalone empty one.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Template Transforms:*") {
        $mark = $d.Range($p.Range.End - 1, $p.Range.End)
        $mark.Delete()
        break
    }
}

# 1d. Fill the still-empty new paragraph with the original "Subjects:" text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Template Transforms:*") {
        $p3 = $p.Next()
        $insPt = $d.Range($p3.Range.Start, $p3.Range.Start)
        $insPt.InsertBefore($subjectsText)
        break
    }
}

# --- 2. Three wording tweaks: "Occurrence : <X>" -> "Occurrence : Statement" ---

$d.Content.Find.Execute( `
    "(Context : SubjectKind, Occurrence : Subject, Attribute : Predicate, Value : Object);", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "(Context : SubjectKind, Occurrence : Statement, Attribute : Predicate, Value : Object);", 2)

$d.Content.Find.Execute( `
    "(Context : PredicateKind, Occurrence : Predicate, Attribute : Subject, Value : Object);", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "(Context : PredicateKind, Occurrence : Statement, Attribute : Subject, Value : Object);", 2)

$d.Content.Find.Execute( `
    "(Context : ObjectKind, Occurrence : Object, Attribute : Subject, Value : Predicate);", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "(Context : ObjectKind, Occurrence : Statement, Attribute : Subject, Value : Predicate);", 2)
